# Update Symptom_Composition_by_ABG_VBG_Overlap sheet data (rows 2-12)
# per "Remove redundant hypercapnia columns" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Administrative",               1.7,  2.5,  1.5),
    @("Diseases (patient-stated)",    5,    7.5,  4.2),
    @("Injuries & adverse effects",   16.4, 22.7, 10.3),
    @("Other",                        3.8,  6.2,  6.6),
    @("Symptom – Circulatory",        6.3,  6.5,  9),
    @("Symptom – Digestive",          7.3,  10.6, 9.9),
    @("Symptom – General",            2.3,  2.9,  4.2),
    @("Symptom – Nervous",            12.2, 9.800000000000001, 13.4),
    @("Symptom – Respiratory",        41.6, 24,   38.5),
    @("Symptom – Skin/Hair/Nails",    0.9,  2.7,  1),
    @("Uncodable/Unknown",            2.5,  4.5,  1.3)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
